$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 686, shifting existing rows 686:766 down to 687:767
$ws.Rows.Item(686).Insert()

# Populate the newly inserted row 686 with the new weekly record
$ws.Range("A686").Value = 10
$ws.Range("B686").Value = "Vega Modelo de Temuco"
$ws.Range("C686").Value = "La Araucanía"
$ws.Range("D686").Value = 45124
$ws.Range("E686").Value = 9
$ws.Range("F686").Value = "Fruta"
$ws.Range("G686").Value = 100108
$ws.Range("H686").Value = "Tropicales y subtropicales"
$ws.Range("I686").Value = 100108005
$ws.Range("J686").Value = "Piña"
$ws.Range("K686").Value = "Caramelo"
$ws.Range("L686").Value = "Segunda"
$ws.Range("M686").Value = 85
$ws.Range("N686").Value = 24000
$ws.Range("O686").Value = 24000
$ws.Range("P686").Value = 24000
$ws.Range("Q686").Value = "$/caja 14 unidades"
$ws.Range("R686").Value = "Ecuador"
$ws.Range("S686").Value = 1714
$ws.Range("T686").Value = 14
